$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (2022-12-15 -> 2022-12-17)
$ws.Name = "Through 2022-12-17"

# Update the column header label for the "through" date
$ws.Range("I1").Value = "2022 (through 12-17)"

# Update November (row 12), December (row 13) and Total (row 14) values in column I
$ws.Range("I12").Value = 119
$ws.Range("I13").Value = 73
$ws.Range("I14").Value = 1590
